$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 7-13: a repeat block (7,gama) + rows 2-6 duplicated + (7,gama) again
$newRows = @(
    @(7, "gama", 37289),
    @(1, "alfa", 36892),
    @(2, "beta", 36892),
    @(3, "gama", 37289),
    @(4, "epsilon", 37683),
    @(5, "theta", 38081),
    @(7, "gama", 37289)
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

# Apply the same date number formatting (style) used by the existing C column
# dates, reusing the existing style record instead of creating a new one.
$ws.Range("C2").Copy()
$ws.Range("C7:C13").PasteSpecial(-4122)

# Update the selection to match the newly added block
$null = $ws.Range("A8:C13").Select()
